$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1308.5741
$ws.Range("J17").Value = 1291.5962
$ws.Range("L17").Value = 3874.7886
$ws.Range("N17").Value = -4210.7886

$ws.Range("H86").Value = 6347.684
$ws.Range("I86").Value = 8734
$ws.Range("J86").Value = 4955.6665
$ws.Range("K86").Value = 8734
$ws.Range("L86").Value = 4955.6665
$ws.Range("M86").Value = -7611
$ws.Range("N86").Value = -7201.6665

$ws.Range("H89").Value = 6347.684
$ws.Range("I89").Value = 8734
$ws.Range("J89").Value = 4955.6665
$ws.Range("K89").Value = 43670
$ws.Range("L89").Value = 24778.3325
$ws.Range("M89").Value = -38054
$ws.Range("N89").Value = -36010.3325

$ws.Range("H116").Value = 7829.8887
$ws.Range("J116").Value = 4611
$ws.Range("L116").Value = 4611
$ws.Range("N116").Value = -11495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 34121.098
$ws.Range("I97").Value = 2039.56
$ws.Range("J97").Value = 167794.17
$ws.Range("K97").Value = 2039.56
$ws.Range("L97").Value = 167794.17
$ws.Range("M97").Value = -1543.56
$ws.Range("N97").Value = -168786.17

$ws.Range("H110").Value = 2016.4103
$ws.Range("I110").Value = 2233.3333
$ws.Range("K110").Value = 2233.3333
$ws.Range("M110").Value = -188.3332999999998

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H132").Value = 2616.8096
$ws.Range("J132").Value = 4000
$ws.Range("L132").Value = 12000
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 64972.2
$ws.Range("J2").Value = 64972.2
$ws.Range("L2").Value = 64972.2
$ws.Range("N2").Value = -65198.2

$ws.Range("H13").Value = 99803
$ws.Range("J13").Value = 99803
$ws.Range("L13").Value = 99803
$ws.Range("N13").Value = -100139

$ws.Range("H107").Value = 2144.5117
$ws.Range("I107").Value = 1843.7878
$ws.Range("K107").Value = 1843.7878
$ws.Range("M107").Value = 76.21219999999994

$ws.Range("H140").Value = 49649.266
$ws.Range("J140").Value = 49649.266
$ws.Range("L140").Value = 49649.266
$ws.Range("N140").Value = -60009.266

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6631.5835
$ws.Range("J58").Value = 6557
$ws.Range("L58").Value = 6557
$ws.Range("N58").Value = -6963

$ws.Range("H99").Value = 11436.657
$ws.Range("I99").Value = 7935.1665
$ws.Range("J99").Value = 15144.117
$ws.Range("K99").Value = 7935.1665
$ws.Range("L99").Value = 15144.117
$ws.Range("M99").Value = -6437.1665
$ws.Range("N99").Value = -18140.117

$ws.Range("H126").Value = 11436.657
$ws.Range("I126").Value = 7935.1665
$ws.Range("J126").Value = 15144.117
$ws.Range("K126").Value = 23805.4995
$ws.Range("L126").Value = 45432.351
$ws.Range("M126").Value = -21335.4995
$ws.Range("N126").Value = -50372.351

$ws.Range("H132").Value = 8302.394
$ws.Range("I132").Value = 2379.3333
$ws.Range("K132").Value = 7137.999899999999
$ws.Range("M132").Value = -4607.999899999999

$ws.Range("H134").Value = 4593.3125
$ws.Range("I134").Value = 4276.1816
$ws.Range("K134").Value = 12828.5448
$ws.Range("M134").Value = -10293.5448

$ws.Range("H136").Value = 6631.5835
$ws.Range("J136").Value = 6557
$ws.Range("L136").Value = 19671
$ws.Range("N136").Value = -24771

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2128.4443
$ws.Range("I5").Value = 1451
$ws.Range("J5").Value = 4499.5
$ws.Range("K5").Value = 4353
$ws.Range("L5").Value = 13498.5
$ws.Range("M5").Value = -4241
$ws.Range("N5").Value = -13722.5

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H88").Value = 9059.546
$ws.Range("I88").Value = 6990.5
$ws.Range("J88").Value = 10241.857
$ws.Range("K88").Value = 20971.5
$ws.Range("L88").Value = 30725.571
$ws.Range("M88").Value = -20543.5
$ws.Range("N88").Value = -31581.571

$ws.Range("H91").Value = 9059.546
$ws.Range("I91").Value = 6990.5
$ws.Range("J91").Value = 10241.857
$ws.Range("K91").Value = 20971.5
$ws.Range("L91").Value = 30725.571
$ws.Range("M91").Value = -19489.5
$ws.Range("N91").Value = -33689.571

$ws.Range("H99").Value = 2179.3333
$ws.Range("I99").Value = 969.6
$ws.Range("K99").Value = 2908.8
$ws.Range("M99").Value = -662.8000000000002

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()

$ws.Range("H131").Value = 1549499.1
$ws.Range("I131").Value = 2263165.5
$ws.Range("K131").Value = 6789496.5
$ws.Range("M131").Value = -6784456.5

$ws.Range("H135").Value = 2128.4443
$ws.Range("I135").Value = 1451
$ws.Range("J135").Value = 4499.5
$ws.Range("K135").Value = 13059
$ws.Range("L135").Value = 40495.5
$ws.Range("M135").Value = -10524
$ws.Range("N135").Value = -45565.5

$ws.Range("H140").Value = 1024.2727
$ws.Range("I140").Value = 1024.2727
$ws.Range("K140").Value = 3072.8181
$ws.Range("M140").Value = 2107.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17369.312
$ws.Range("I70").Value = 23674.7
$ws.Range("K70").Value = 23674.7
$ws.Range("M70").Value = -23404.7

$ws.Range("H73").Value = 17369.312
$ws.Range("I73").Value = 23674.7
$ws.Range("K73").Value = 23674.7
$ws.Range("M73").Value = -22738.7

$ws.Range("H80").Value = 25560332
$ws.Range("I80").Value = 46004200
$ws.Range("J80").Value = 5499.5
$ws.Range("K80").Value = 46004200
$ws.Range("L80").Value = 5499.5
$ws.Range("M80").Value = -46003202
$ws.Range("N80").Value = -7495.5

$ws.Range("H83").Value = 25560332
$ws.Range("I83").Value = 46004200
$ws.Range("J83").Value = 5499.5
$ws.Range("K83").Value = 230021000
$ws.Range("L83").Value = 27497.5
$ws.Range("M83").Value = -230016008
$ws.Range("N83").Value = -37481.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 100250
$ws.Range("J6").Value = 100250
$ws.Range("L6").Value = 100250
$ws.Range("N6").Value = -100474

$ws.Range("H132").Value = 19542.158
$ws.Range("I132").Value = 25119.521
$ws.Range("K132").Value = 75358.56299999999
$ws.Range("M132").Value = -72828.56299999999

$ws.Range("H136").Value = 9003.25
$ws.Range("I136").Value = 11747.056
$ws.Range("K136").Value = 35241.16800000001
$ws.Range("M136").Value = -32691.16800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1513.0714
$ws.Range("I100").Value = 1298.2
$ws.Range("J100").Value = 2050.25
$ws.Range("K100").Value = 2596.4
$ws.Range("L100").Value = 4100.5
$ws.Range("M100").Value = -2055.4
$ws.Range("N100").Value = -5182.5

$ws.Range("H101").Value = 50174.855
$ws.Range("J101").Value = 50174.855
$ws.Range("L101").Value = 50174.855
$ws.Range("N101").Value = -56664.855

$ws.Range("H118").Value = 106666.336
$ws.Range("J118").Value = 106666.336
$ws.Range("L118").Value = 106666.336
$ws.Range("N118").Value = -109980.336

$ws.Range("H132").Value = 6997.7144
$ws.Range("I132").Value = 6929.4653
$ws.Range("J132").Value = 7486.8335
$ws.Range("K132").Value = 20788.3959
$ws.Range("L132").Value = 22460.5005
$ws.Range("M132").Value = -18258.3959
$ws.Range("N132").Value = -27520.5005
